# BRBW_YR_FIN.xlsx — "Doing Updates for Financials"
#
# The three "Period Ending" header rows (7, 38, 80) get a brand-new most-recent
# fiscal year inserted before column D (old D..J shift right into E..K, and a
# new date — 2017-12-31, serial 43100 — lands in D).
#
# Every other data row keeps its D..J values exactly where they are; only a
# new trailing column K is populated, and its value is simply a copy of
# column J (the last existing year) for that row. A brand-new, still-empty
# column L is appended after K everywhere K already exists, carrying over K's
# (pre-existing, already-correct) number format/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Stamp column L into existence by copying the (currently empty, but
#    already correctly styled) column K cells one column to the right.
#    This must happen BEFORE we touch K's values below, since K is still
#    blank everywhere at this point (mirrors the source workbook state) and
#    this picks up each row's existing style (date style for row 7/38/80,
#    numeric style for the data rows) without fabricating any new style.
# ---------------------------------------------------------------------
$ws.Range("K7:K35").Copy($ws.Range("L7:L35"))
$ws.Range("K38:K77").Copy($ws.Range("L38:L77"))
$ws.Range("K80:K102").Copy($ws.Range("L80:L102"))

# ---------------------------------------------------------------------
# 2) The three "Period Ending" date header rows: shift D..J right into
#    E..K, then drop the new fiscal year end date (2017-12-31) into D.
#    Walk right-to-left so each source cell is read before it is
#    overwritten.
# ---------------------------------------------------------------------
$headerRows = @(7, 38, 80)
foreach ($r in $headerRows) {
    $ws.Range("K$r").Value = $ws.Range("J$r").Value2
    $ws.Range("J$r").Value = $ws.Range("I$r").Value2
    $ws.Range("I$r").Value = $ws.Range("H$r").Value2
    $ws.Range("H$r").Value = $ws.Range("G$r").Value2
    $ws.Range("G$r").Value = $ws.Range("F$r").Value2
    $ws.Range("F$r").Value = $ws.Range("E$r").Value2
    $ws.Range("E$r").Value = $ws.Range("D$r").Value2
    $ws.Range("D$r").Value = 43100
}

# ---------------------------------------------------------------------
# 3) Every other populated data row: column K simply repeats column J's
#    value (blank "section header" rows such as 11/16/19/39/40/55/56/67/
#    82/90/95 stay blank and need no value here).
# ---------------------------------------------------------------------
$dataRows = @(
    8,9,10,12,13,14,15,17,18,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,
    41,42,43,44,45,46,47,48,49,50,51,52,53,54,
    57,58,59,60,61,62,63,64,65,66,68,69,70,71,72,73,74,75,76,77,
    81,83,84,85,86,87,88,89,91,92,93,94,96,97,98,99,100,101,102
)
foreach ($r in $dataRows) {
    $ws.Range("K$r").Value = $ws.Range("J$r").Value2
}
